$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# DownloadEmailAttachments option: No -> Yes
$ws.Range("B4").Value = "Yes"

# Move active selection (cosmetic, matches saved cursor position)
$ws.Range("B7").Select()
